$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data cells: row 41 (Frágangur) week "Vika 4"
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 0.5

# Recalculate formulas (SUM totals depend on C41/D41)
$excel.Calculate()

# Update window / view state to match the authored selection
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("D43").Select()

$wb.Windows.Item(1).WindowState = -4143
